# Resort the worksheet tabs: move "总计" (the summary/total sheet) so it
# becomes the first sheet, ahead of "2022-Q2". The underlying worksheet
# contents travel with their sheet (COM Move keeps name/content bound
# together), so no cell data is rewritten here - only the tab order changes.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$firstSheet = $wb.Worksheets.Item(1)

$totalSheet.Move($firstSheet)
